$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "26-09-2025"
$ws.Range("B11").Value = "The price of gold in India today is ₹11,488 per gram for 24 karat gold, ₹10,530 per gram for 22 karat gold and ₹8,616 per gram for 18 karat gold (also called 999 gold)."
